# Refresh crypto price/volume data (GitHub Actions run Mon Jun 17 03:29:32 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that were authored as *text* (some use "."
# as a thousands separator, e.g. "66.345.75"). Assigning a plain numeric-
# looking string via .Value auto-converts it to a real number, which does
# not match the source data. Force text storage (without leaving a lasting
# NumberFormat change) by flipping the cell to Text format just for the
# write, then restoring the default "Normal" style.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "66.345.75"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.587.64"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue "D5" "606.82"
$ws.Range("E5").Value = "  +0.18%  "
Set-TextValue "D6" "148.09"
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("D7").Value = "3.587.48"
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.86%  "
Set-TextValue "D10" "0.135"
$ws.Range("E10").Value = "  -0.30%  "
Set-TextValue "D11" "7.82"
$ws.Range("E11").Value = "  +0.03%  "
Set-TextValue "D12" "0.413"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "4.197.33"
$ws.Range("E13").Value = "  +0.86%  "
Set-TextValue "D14" "0.0000205"
$ws.Range("E14").Value = "  -0.68%  "
Set-TextValue "D15" "29.48"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "3.577.37"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("D18").Value = "66.405.45"
$ws.Range("E18").Value = "  +0.29%  "
Set-TextValue "D19" "11.04"
$ws.Range("E19").Value = "  -3.04%  "
Set-TextValue "D20" "6.31"
$ws.Range("E20").Value = "  +1.86%  "
Set-TextValue "D21" "14.85"
$ws.Range("E21").Value = "  +1.10%  "
Set-TextValue "D22" "422.66"
$ws.Range("E22").Value = "  -1.89%  "
Set-TextValue "D23" "0.610"
$ws.Range("E23").Value = "  +0.12%  "
Set-TextValue "D24" "78.53"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("E25").Value = "  +0.02%  "
Set-TextValue "D26" "0.0000120"
$ws.Range("E26").Value = "  +2.39%  "
Set-TextValue "D27" "8.19"
$ws.Range("E27").Value = "  +4.33%  "
Set-TextValue "D28" "9.34"
$ws.Range("E28").Value = "  +2.63%  "
Set-TextValue "D29" "2.49"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "3.585.14"
$ws.Range("E31").Value = "  +0.84%  "
Set-TextValue "D32" "0.157"
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "25.04"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D34" "1.43"
$ws.Range("E34").Value = "  -1.06%  "
Set-TextValue "D36" "7.73"
$ws.Range("E36").Value = "  -0.85%  "
Set-TextValue "D37" "5.56"
$ws.Range("E37").Value = "  +0.33%  "
Set-TextValue "D38" "1.66"
$ws.Range("E38").Value = "  -2.96%  "
Set-TextValue "D39" "175.02"
$ws.Range("E39").Value = "  +0.86%  "
Set-TextValue "D40" "0.0850"
$ws.Range("E40").Value = "  +0.49%  "
Set-TextValue "D41" "5.17"
$ws.Range("E41").Value = "  -0.06%  "
Set-TextValue "D42" "0.881"
$ws.Range("E42").Value = "  -1.32%  "
Set-TextValue "D43" "45.84"
$ws.Range("E43").Value = "  -0.41%  "
Set-TextValue "D44" "1.85"
$ws.Range("E44").Value = "  -4.14%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +4.59%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D47" "7.14"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "23.49"
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D49" "24.07"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("E50").Value = "  -5.87%  "
Set-TextValue "D51" "0.951"
$ws.Range("E51").Value = "  +2.18%  "
